$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "ÍNDICE") {
        continue
    }

    # Turn what used to be the first data row (numeric sounding/volume pair)
    # into a second, repeated header row - same text & style as row 5.
    $ws.Range("A6").Value = "Sondagem (mm)"
    $ws.Range("B6").Value = "Volume (litros)"

    $ws.Range("A5:B5").Copy() | Out-Null
    $ws.Range("A6:B6").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false
